$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw-data")

# The single "date" column (G) is split into three columns: date.year,
# date.month, date.day. Insert two extra columns right after G so the
# existing date values (G) can be decomposed into G/H/I.
$ws.Columns("H:I").Insert()

$lastRow = 31
for ($r = 2; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 7).Value2
    if ($d -ne $null -and $d -ne "") {
        $dt = [DateTime]::FromOADate($d)
        $ws.Cells.Item($r, 7).Value2 = $dt.Year
        $ws.Cells.Item($r, 8).Value2 = $dt.Month
        $ws.Cells.Item($r, 9).Value2 = $dt.Day
    }
}

# These columns no longer hold dates, so drop the date number format
# (ClearFormats reverts to the default/general style instead of stamping a
# brand new number format onto the style table).
$ws.Columns("G:I").ClearFormats() | Out-Null

# Header row for the new columns.
$ws.Cells.Item(1,7).Value2 = "date.year"
$ws.Cells.Item(1,8).Value2 = "date.month"
$ws.Cells.Item(1,9).Value2 = "date.day"

# Keep the autofilter / used range in sync with the new column count.
$ws.AutoFilterMode = $false
$ws.Range("A1:L$lastRow").AutoFilter() | Out-Null

# _xlnm._FilterDatabase still points at the old J column; refresh it.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='raw-data'!`$A`$1:`$L`$$lastRow"
    }
}

# The edit session left the cursor on H6 (month column) in the original file.
$ws.Range("H6").Select() | Out-Null

Write-Host "done"
